# Auto-generated edit script: update Sheets via scheduled runner
# Applies cell-value corrections to the Leve profit tables across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1780
$ws.Range("J32").Value = 2300
$ws.Range("L32").Value = 2300
$ws.Range("N32").Value = -2952
$ws.Range("H64").Value = 3952.5789
$ws.Range("I64").Value = 2985.5715
$ws.Range("J64").Value = 4516.6665
$ws.Range("K64").Value = 2985.5715
$ws.Range("L64").Value = 4516.6665
$ws.Range("M64").Value = -2737.5715
$ws.Range("N64").Value = -5012.6665
$ws.Range("H67").Value = 3952.5789
$ws.Range("I67").Value = 2985.5715
$ws.Range("J67").Value = 4516.6665
$ws.Range("K67").Value = 2985.5715
$ws.Range("L67").Value = 4516.6665
$ws.Range("M67").Value = -2127.5715
$ws.Range("N67").Value = -6232.6665
$ws.Range("H96").Value = 947
$ws.Range("I96").Value = 1005.2
$ws.Range("K96").Value = 3015.6
$ws.Range("M96").Value = -1642.6
$ws.Range("H97").Value = 2194.5557
$ws.Range("J97").Value = 2194.5557
$ws.Range("L97").Value = 6583.6671
$ws.Range("N97").Value = -7575.6671
$ws.Range("H107").Value = 722.0769
$ws.Range("I107").Value = 1041
$ws.Range("J107").Value = 522.75
$ws.Range("K107").Value = 1041
$ws.Range("L107").Value = 522.75
$ws.Range("M107").Value = 879
$ws.Range("N107").Value = -4362.75
$ws.Range("H111").Value = 4147
$ws.Range("I111").Value = 4147
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 12441
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -9374
$ws.Range("H116").Value = 4233.8335
$ws.Range("I116").Value = 1598.75
$ws.Range("J116").Value = 5551.375
$ws.Range("K116").Value = 1598.75
$ws.Range("L116").Value = 5551.375
$ws.Range("M116").Value = 1843.25
$ws.Range("N116").Value = -12435.375
$ws.Range("H137").Value = 52387.9
$ws.Range("I137").Value = 2550.4666
$ws.Range("J137").Value = 201900.2
$ws.Range("K137").Value = 7651.399800000001
$ws.Range("L137").Value = 605700.6000000001
$ws.Range("M137").Value = -5101.399800000001
$ws.Range("N137").Value = -610800.6000000001
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1303.6451
$ws.Range("I2").Value = 1289.5264
$ws.Range("J2").Value = 1326
$ws.Range("K2").Value = 1289.5264
$ws.Range("L2").Value = 1326
$ws.Range("M2").Value = -1176.5264
$ws.Range("N2").Value = -1552
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 12
$ws.Range("H45").Value = 2661.3333
$ws.Range("I45").Value = 3132.2222
$ws.Range("K45").Value = 3132.2222
$ws.Range("M45").Value = -2755.2222
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("H74").Value = 166674220
$ws.Range("I74").Value = 333347000
$ws.Range("J74").Value = 1438
$ws.Range("K74").Value = 333347000
$ws.Range("L74").Value = 1438
$ws.Range("M74").Value = -333346126
$ws.Range("N74").Value = -3186
$ws.Range("H77").Value = 166674220
$ws.Range("I77").Value = 333347000
$ws.Range("J77").Value = 1438
$ws.Range("K77").Value = 1666735000
$ws.Range("L77").Value = 7190
$ws.Range("M77").Value = -1666730632
$ws.Range("N77").Value = -15926
$ws.Range("H105").Value = 38141.668
$ws.Range("J105").Value = 38141.668
$ws.Range("L105").Value = 38141.668
$ws.Range("N105").Value = -45129.668
$ws.Range("H116").Value = 1303.6451
$ws.Range("I116").Value = 1289.5264
$ws.Range("J116").Value = 1326
$ws.Range("K116").Value = 1289.5264
$ws.Range("L116").Value = 1326
$ws.Range("M116").Value = 1004.4736
$ws.Range("N116").Value = -5914
$ws.Range("N59").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1303.6451
$ws.Range("I3").Value = 1289.5264
$ws.Range("J3").Value = 1326
$ws.Range("K3").Value = 1289.5264
$ws.Range("L3").Value = 1326
$ws.Range("M3").Value = -1175.5264
$ws.Range("N3").Value = -1554
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 15
$ws.Range("H11").Value = 231
$ws.Range("I11").Value = 231
$ws.Range("K11").Value = 231
$ws.Range("M11").Value = -91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8066.1914
$ws.Range("I31").Value = 10304.906
$ws.Range("J31").Value = 3290.2666
$ws.Range("K31").Value = 10304.906
$ws.Range("L31").Value = 3290.2666
$ws.Range("M31").Value = -10009.906
$ws.Range("N31").Value = -3880.2666
$ws.Range("H34").Value = 8066.1914
$ws.Range("I34").Value = 10304.906
$ws.Range("J34").Value = 3290.2666
$ws.Range("K34").Value = 10304.906
$ws.Range("L34").Value = 3290.2666
$ws.Range("M34").Value = -10102.906
$ws.Range("N34").Value = -3694.2666
$ws.Range("H41").Value = 30000
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30856
$ws.Range("H47").Value = 18000
$ws.Range("J47").Value = 18000
$ws.Range("L47").Value = 18000
$ws.Range("N47").Value = -19132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 101895.555
$ws.Range("I36").Value = 1884
$ws.Range("K36").Value = 5652
$ws.Range("M36").Value = -5483
$ws.Range("H112").Value = 2333.25
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 3166.5
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 9499.5
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -11715.5
$ws.Range("H131").Value = 735.39
$ws.Range("J131").Value = 748.29785
$ws.Range("L131").Value = 2244.89355
$ws.Range("N131").Value = -12324.89355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 12504610
$ws.Range("I58").Value = 2293.3333
$ws.Range("K58").Value = 2293.3333
$ws.Range("M58").Value = -2016.3333
$ws.Range("H132").Value = 99915.44
$ws.Range("I132").Value = 90095.75
$ws.Range("J132").Value = 129374.5
$ws.Range("K132").Value = 270287.25
$ws.Range("L132").Value = 388123.5
$ws.Range("M132").Value = -267757.25
$ws.Range("N132").Value = -393183.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8415.5
$ws.Range("I7").Value = 8681.875
$ws.Range("J7").Value = 7350
$ws.Range("K7").Value = 8681.875
$ws.Range("L7").Value = 7350
$ws.Range("M7").Value = -8569.875
$ws.Range("N7").Value = -7574
$ws.Range("H55").Value = 76.210526
$ws.Range("I55").Value = 39.1
$ws.Range("J55").Value = 117.44444
$ws.Range("K55").Value = 39.1
$ws.Range("L55").Value = 117.44444
$ws.Range("M55").Value = 133.9
$ws.Range("N55").Value = -463.44444
$ws.Range("H57").Value = 9490
$ws.Range("H126").Value = 8415.5
$ws.Range("I126").Value = 8681.875
$ws.Range("J126").Value = 7350
$ws.Range("K126").Value = 26045.625
$ws.Range("L126").Value = 22050
$ws.Range("M126").Value = -23575.625
$ws.Range("N126").Value = -26990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1905.6
$ws.Range("I132").Value = 1087.625
$ws.Range("J132").Value = 2840.4285
$ws.Range("K132").Value = 3262.875
$ws.Range("L132").Value = 8521.2855
$ws.Range("N132").Value = -13581.2855
$ws.Range("M132").Value = -732.875
